$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") for rows 2-11 from serial date 45212 to 45221.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = 45221
}
